$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing the "LOM3206 - Eletrônica (Requisito)" requirement (row 23),
# shifting the subsequent requirement rows up.
$ws.Rows.Item(23).Delete()

# Update remaining requirement rows text (now shifted up one row).
$ws.Range("B24").Value = "LOM3234 -  Óptica Física  (Requisito)`n"
$ws.Range("C24").Value = "LOM3234 -  Óptica Física  (Requisito)`n"

$ws.Range("B25").Value = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"
$ws.Range("C25").Value = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)`n"

$wb.Save()
